# "Generate Report for Handoff"
#
# The a84af520-72cf-4773-8579-da167f94a862 file is no longer tracked, so its
# row (row 3) is removed from every sheet (Overview, zh-cn, de-de), shifting
# the ".localization-config" row up from row 4 to row 3. The outstanding
# "Handed back: in sync with en-US" status for 8aae9dc6-... becomes
# "Ready for handoff", and its Latest Handoff Datetime is refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Drop the a84af520-...-md row (row 3) from every sheet -----------------
# This shifts the ".localization-config" row up from row 4 to row 3 and
# keeps cell styling / shared values intact.
$overview.Rows(3).Delete()
$zhcn.Rows(3).Delete()
$dede.Rows(3).Delete()

# --- Refresh the handoff status text ---------------------------------------
$overview.Cells.Replace("Handed back: in sync with en-US", "Ready for handoff")
$zhcn.Cells.Replace("Handed back: in sync with en-US", "Ready for handoff")
$dede.Cells.Replace("Handed back: in sync with en-US", "Ready for handoff")

# --- Refresh the Latest Handoff Datetime for the new handoff ---------------
$zhcn.Cells.Replace("2016-02-06 03:59:21", "2016-02-06 04:01:03")
$dede.Cells.Replace("2016-02-06 03:59:32", "2016-02-06 04:01:14")

# --- Rebuild hyperlinks (row delete does not renumber/relocate them) -------
$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ba27dad74b53a60d7dd70825ac329f0a189d788b/e2e/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md")
$overview.Hyperlinks.Add($overview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ba27dad74b53a60d7dd70825ac329f0a189d788b/.localization-config", $null, $null, ".localization-config")

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ba27dad74b53a60d7dd70825ac329f0a189d788b/e2e/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1112658416616796686452b0f2ec70ae2eb5b12a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.zh-cn.xlf", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9adc6fa52fab058cb0357755e10446a209a79183/e2e/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/92abf0d33a47fbac90e821b3b22bfbf7c3741081/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.zh-cn.xlf", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ba27dad74b53a60d7dd70825ac329f0a189d788b/.localization-config", $null, $null, ".localization-config")

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ba27dad74b53a60d7dd70825ac329f0a189d788b/e2e/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md")
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aceb382cc1caf259cee746943df7206fe5aa1cdc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.de-de.xlf", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d5d7613cd67ad4ecbe0449b9ff0be8b94ce0002b/e2e/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/314028db5985e8488a4bd21408c107808618889c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.de-de.xlf", $null, $null, "8aae9dc6-d75c-4ce4-99c9-e106ad496bd8.d0d7122f89f9e7b41583b5e8da0c92c9b295a988.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ba27dad74b53a60d7dd70825ac329f0a189d788b/.localization-config", $null, $null, ".localization-config")
